$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 444, shifting existing rows 444:468 down to 445:469.
$ws.Rows.Item(444).Insert()

# Populate the new row 444 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,J,K,Q,T are identical to the surrounding rows
# for this market/product subset; only D,L,M,N,O,P,R,S differ per record.
$ws.Range("A444").Value = 7
$ws.Range("B444").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C444").Value = "Ñuble"
$ws.Range("D444").Value = 44516
$ws.Range("E444").Value = 16
$ws.Range("F444").Value = "Fruta"
$ws.Range("G444").Value = 100102
$ws.Range("H444").Value = "Cítricos"
$ws.Range("I444").Value = 100102003
$ws.Range("J444").Value = "Limón"
$ws.Range("K444").Value = "Sin especificar"
$ws.Range("L444").Value = "1a amarillo"
$ws.Range("M444").Value = 160
$ws.Range("N444").Value = 7000
$ws.Range("O444").Value = 7500
$ws.Range("P444").Value = 7250
$ws.Range("Q444").Value = "$/malla 16 kilos"
$ws.Range("R444").Value = "Región de O'Higgins"
$ws.Range("S444").Value = 453
$ws.Range("T444").Value = 16
